$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1716.5625
$ws.Range("I43").Value = 1523.6364
$ws.Range("J43").Value = 2141
$ws.Range("K43").Value = 1523.6364
$ws.Range("L43").Value = 2141
$ws.Range("M43").Value = -1454.6364
$ws.Range("N43").Value = -2279
$ws.Range("H58").Value = 1419.8334
$ws.Range("J58").Value = 3499.5
$ws.Range("L58").Value = 10498.5
$ws.Range("N58").Value = -10798.5
$ws.Range("H61").Value = 142857520
$ws.Range("I61").Value = 142857520
$ws.Range("K61").Value = 428572560
$ws.Range("M61").Value = -428572388
$ws.Range("H62").Value = 9279
$ws.Range("I62").Value = 11917.333
$ws.Range("K62").Value = 11917.333
$ws.Range("M62").Value = -11293.333
$ws.Range("H65").Value = 9279
$ws.Range("I65").Value = 11917.333
$ws.Range("K65").Value = 59586.665
$ws.Range("M65").Value = -56466.665
$ws.Range("H137").Value = 108898.7
$ws.Range("I137").Value = 256211.58
$ws.Range("K137").Value = 768634.74
$ws.Range("M137").Value = -766084.74
$ws.Range("H138").Value = 5556.8037
$ws.Range("I138").Value = 4494.5
$ws.Range("J138").Value = 5647.213
$ws.Range("K138").Value = 13483.5
$ws.Range("L138").Value = 16941.639
$ws.Range("M138").Value = -8343.5
$ws.Range("N138").Value = -27221.639
$ws.Range("H141").Value = 7215.0435
$ws.Range("I141").Value = 7748.15
$ws.Range("K141").Value = 23244.45
$ws.Range("M141").Value = -18064.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 1295.1428
$ws.Range("I31").Value = 1295.1428
$ws.Range("K31").Value = 1295.1428
$ws.Range("M31").Value = -1001.1428
$ws.Range("H32").Value = 10924.5
$ws.Range("I32").Value = 9603.012000000001
$ws.Range("K32").Value = 9603.012000000001
$ws.Range("M32").Value = -9316.012000000001
$ws.Range("H45").Value = 9527324
$ws.Range("I45").Value = 17858240
$ws.Range("J45").Value = 6277
$ws.Range("K45").Value = 17858240
$ws.Range("L45").Value = 6277
$ws.Range("M45").Value = -17857863
$ws.Range("N45").Value = -7031
$ws.Range("H61").Value = 5132.037
$ws.Range("I61").Value = 5271.7915
$ws.Range("K61").Value = 5271.7915
$ws.Range("M61").Value = -5059.7915
$ws.Range("H74").Value = 56636.293
$ws.Range("I74").Value = 4620.5625
$ws.Range("K74").Value = 4620.5625
$ws.Range("M74").Value = -3746.5625
$ws.Range("H77").Value = 56636.293
$ws.Range("I77").Value = 4620.5625
$ws.Range("K77").Value = 23102.8125
$ws.Range("M77").Value = -18734.8125
$ws.Range("H80").Value = 29996.5
$ws.Range("I80").Value = 29996.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 29996.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -28998.5
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 29996.5
$ws.Range("I83").Value = 29996.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 89989.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -84997.5
$ws.Range("N83").ClearContents()
$ws.Range("H88").Value = 4192.364
$ws.Range("J88").Value = 3116.5
$ws.Range("L88").Value = 3116.5
$ws.Range("N88").Value = -3928.5
$ws.Range("H91").Value = 4192.364
$ws.Range("J91").Value = 3116.5
$ws.Range("L91").Value = 3116.5
$ws.Range("N91").Value = -5924.5
$ws.Range("H136").Value = 5132.037
$ws.Range("I136").Value = 5271.7915
$ws.Range("K136").Value = 15815.3745
$ws.Range("M136").Value = -13265.3745
$ws.Range("H139").Value = 96714.664
$ws.Range("J139").Value = 96714.664
$ws.Range("L139").Value = 96714.664
$ws.Range("N139").Value = -106994.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I86").Value = 36114252
$ws.Range("J86").Value = 5012.636
$ws.Range("K86").Value = 36114252
$ws.Range("L86").Value = 5012.636
$ws.Range("M86").Value = -36113129
$ws.Range("N86").Value = -7258.636
$ws.Range("I89").Value = 36114252
$ws.Range("J89").Value = 5012.636
$ws.Range("K89").Value = 180571260
$ws.Range("L89").Value = 25063.18
$ws.Range("M89").Value = -180565644
$ws.Range("N89").Value = -36295.18
$ws.Range("H94").Value = 2846480.5
$ws.Range("I94").Value = 3789240
$ws.Range("K94").Value = 3789240
$ws.Range("M94").Value = -3788789
$ws.Range("H134").Value = 24598.8
$ws.Range("I134").Value = 3497
$ws.Range("K134").Value = 10491
$ws.Range("M134").Value = -7956
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 842.7778
$ws.Range("I16").Value = 519.9286
$ws.Range("K16").Value = 519.9286
$ws.Range("M16").Value = -232.9286
$ws.Range("H31").Value = 37681.105
$ws.Range("I31").Value = 3474.2856
$ws.Range("J31").Value = 71887.92999999999
$ws.Range("K31").Value = 3474.2856
$ws.Range("L31").Value = 71887.92999999999
$ws.Range("M31").Value = -3179.2856
$ws.Range("N31").Value = -72477.92999999999
$ws.Range("H34").Value = 37681.105
$ws.Range("I34").Value = 3474.2856
$ws.Range("J34").Value = 71887.92999999999
$ws.Range("K34").Value = 3474.2856
$ws.Range("L34").Value = 71887.92999999999
$ws.Range("M34").Value = -3272.2856
$ws.Range("N34").Value = -72291.92999999999
$ws.Range("H69").Value = 26064.6
$ws.Range("I69").Value = 24593.25
$ws.Range("K69").Value = 24593.25
$ws.Range("M69").Value = -23844.25
$ws.Range("H72").Value = 26064.6
$ws.Range("I72").Value = 24593.25
$ws.Range("K72").Value = 73779.75
$ws.Range("M72").Value = -70035.75
$ws.Range("H99").Value = 4449.8237
$ws.Range("J99").Value = 5392.7144
$ws.Range("L99").Value = 5392.7144
$ws.Range("N99").Value = -8388.714400000001
$ws.Range("H113").Value = 842.7778
$ws.Range("I113").Value = 519.9286
$ws.Range("K113").Value = 519.9286
$ws.Range("M113").Value = 1650.0714
$ws.Range("H126").Value = 4449.8237
$ws.Range("J126").Value = 5392.7144
$ws.Range("L126").Value = 16178.1432
$ws.Range("N126").Value = -21118.1432
$ws.Range("H132").Value = 72149.7
$ws.Range("I132").Value = 57412.945
$ws.Range("K132").Value = 172238.835
$ws.Range("M132").Value = -169708.835
$ws.Range("H134").Value = 2293.5
$ws.Range("I134").Value = 1634.7693
$ws.Range("K134").Value = 4904.3079
$ws.Range("M134").Value = -2369.3079
$ws.Range("H141").Value = 553999.9
$ws.Range("J141").Value = 610000
$ws.Range("L141").Value = 610000
$ws.Range("N141").Value = -620360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2795.4285
$ws.Range("J113").Value = 3113.6
$ws.Range("L113").Value = 9340.799999999999
$ws.Range("N113").Value = -13680.8
$ws.Range("H131").Value = 13898683
$ws.Range("I131").Value = 9261949
$ws.Range("J131").Value = 15885855
$ws.Range("K131").Value = 27785847
$ws.Range("L131").Value = 47657565
$ws.Range("M131").Value = -27780807
$ws.Range("N131").Value = -47667645
$ws.Range("H137").Value = 3430.524
$ws.Range("J137").Value = 4001.7856
$ws.Range("L137").Value = 12005.3568
$ws.Range("N137").Value = -22205.3568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2577469.5
$ws.Range("I80").Value = 5690207.5
$ws.Range("K80").Value = 5690207.5
$ws.Range("M80").Value = -5689209.5
$ws.Range("H83").Value = 2577469.5
$ws.Range("I83").Value = 5690207.5
$ws.Range("K83").Value = 28451037.5
$ws.Range("M83").Value = -28446045.5
$ws.Range("H122").Value = 164389.31
$ws.Range("I122").Value = 191530.4
$ws.Range("J122").Value = 4935.375
$ws.Range("K122").Value = 574591.2
$ws.Range("L122").Value = 14806.125
$ws.Range("M122").Value = -572141.2
$ws.Range("N122").Value = -19706.125
$ws.Range("H132").Value = 2705.2
$ws.Range("I132").Value = 2074.318
$ws.Range("K132").Value = 6222.954000000001
$ws.Range("M132").Value = -3692.954000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1699333.6
$ws.Range("J2").Value = 39200.4
$ws.Range("L2").Value = 39200.4
$ws.Range("N2").Value = -39424.4
$ws.Range("H68").Value = 5400.75
$ws.Range("I68").Value = 2601
$ws.Range("K68").Value = 2601
$ws.Range("M68").Value = -1852
$ws.Range("H71").Value = 5400.75
$ws.Range("I71").Value = 2601
$ws.Range("K71").Value = 13005
$ws.Range("M71").Value = -9261
$ws.Range("H136").Value = 146336.86
$ws.Range("I136").Value = 169643.08
$ws.Range("J136").Value = 6499.5
$ws.Range("K136").Value = 508929.24
$ws.Range("L136").Value = 19498.5
$ws.Range("M136").Value = -506379.24
$ws.Range("N136").Value = -24598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3073.8386
$ws.Range("I122").Value = 1761.5652
$ws.Range("J122").Value = 6846.625
$ws.Range("K122").Value = 5284.6956
$ws.Range("L122").Value = 20539.875
$ws.Range("M122").Value = -2834.6956
$ws.Range("N122").Value = -25439.875
$ws.Range("H132").Value = 56163804
$ws.Range("I132").Value = 76925940
$ws.Range("K132").Value = 230777820
$ws.Range("M132").Value = -230775290
